$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.411.44'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '2.316.12'
$ws.Range('E3').Value = '  +1.09%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '518.17'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.11'
$ws.Range('D6').NumberFormat = "General"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.80%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').Value = '2.335.93'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.102'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.75%  '
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('E12').Value = '  +3.65%  '
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.79'
$ws.Range('D14').NumberFormat = "General"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').Value = '2.734.21'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '56.525.18'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '2.336.54'
$ws.Range('E18').Value = '  +2.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.43'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.22'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '322.84'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.26%  '
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').NumberFormat = "General"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.75'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').Value = '  +5.47%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.994'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('E27').Value = '  +5.95%  '
$ws.Range('E28').Value = '  +11.34%  '
$ws.Range('E29').Value = '  +3.60%  '
$ws.Range('E30').Value = '  +3.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '166.29'
$ws.Range('D31').NumberFormat = "General"
$ws.Range('D31').Style = "Normal"
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.33'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.994'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('E38').Value = '  +3.57%  '
$ws.Range('E39').Value = '  +6.19%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '37.90'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.381'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '139.21'
$ws.Range('D42').NumberFormat = "General"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('E43').Value = '  +4.45%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.25'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.73%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '279.95'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +7.33%  '
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0506'
$ws.Range('D47').NumberFormat = "General"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('E48').Value = '  +1.84%  '
$ws.Range('E49').Value = '  +2.46%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.78'
$ws.Range('D51').NumberFormat = "General"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +8.05%  '
